$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = @(2.093638017770559, 1.988258088193845, 1.924823135559052, 1.899290298110429, 1.895069714919657, 1.92447750768639, 2.057038668061864, 2.327147835593223, 2.531938890989068, 2.626515185610572, 2.662534732410336, 2.654768104527591, 2.629474404312134, 2.614008111285102, 2.525786764280895, 2.472029831258737, 2.441243408019318, 2.430842462478722, 2.477738554710868, 2.636898182540904, 2.74211728923126, 2.685849492774025, 2.475157269848921, 2.252973449046294)
for ($i = 0; $i -lt $colB.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}

$colC = @(0.315208144334747, 0.2801567759861427, 0.2587678202824577, 0.2500840293490967, 0.2486440194843453, 0.2586505778245396, 0.3030942743851028, 0.3913507429581387, 0.4569399331048771, 0.4869571536984836, 0.4983510753007181, 0.4958959693657903, 0.4878939898856629, 0.4829961026941305, 0.4549819615482988, 0.4378431781427707, 0.4280023127554955, 0.4246732353576022, 0.4396658728004468, 0.4902436220904747, 0.52345731194157, 0.5057156919546628, 0.4388417938517364, 0.3673492968682694)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}

$colE = @(0.03627480459440591, 0.0355052710502779, 0.03502449253365114, 0.03482645905580561, 0.0347934473017677, 0.03502183037607587, 0.03601117015721478, 0.03788696263893776, 0.0392279631514203, 0.03983040632416923, 0.0400574803379854, 0.04000862242515524, 0.03984910885886528, 0.03975126536506757, 0.03918844222106443, 0.03884124884629081, 0.03864083652308281, 0.03857285641741015, 0.03887828202998733, 0.03989599031634583, 0.04055496925413671, 0.04020381186888322, 0.03886154184122148, 0.03738618378678105)
for ($i = 0; $i -lt $colE.Length; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}

$colF = @(0.4443680307746263, 0.387822817061874, 0.3531389305168915, 0.3390132514313251, 0.336668177824194, 0.3529483938344953, 0.4248636149813336, 0.5661985755041457, 0.6702781546542269, 0.7176906081379002, 0.7356546913071611, 0.7317853510981394, 0.7191683204515869, 0.7114413442032514, 0.6671810134426437, 0.6400460337215605, 0.6244449056556647, 0.6191636801734006, 0.6429339538360921, 0.7228739723492197, 0.7751780083420101, 0.7472568307916134, 0.6416283278902171, 0.5279251897347308)
for ($i = 0; $i -lt $colF.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}

$colG = @(0.002542814170109375, 0.002548453507922864, 0.002552096365496674, 0.002553626352529557, 0.002553883158156492, 0.002552116815120353, 0.002544721292929434, 0.002531641837716735, 0.002522889691963589, 0.002519092082156152, 0.002517680286336217, 0.002517983175662045, 0.002518975407081975, 0.002519586594955873, 0.002523141559580103, 0.002525369379284448, 0.002526668070129824, 0.002527110761626039, 0.002525130434088336, 0.002518683252679846, 0.002514622737049388, 0.00251677595289978, 0.002525238405399545, 0.002535028882079347)
for ($i = 0; $i -lt $colG.Length; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $colG[$i]
}

$colI = @(1.646081954904304, 1.6417606853405, 1.640033185484903, 1.639561057955532, 1.639496636650911, 1.640025880760099, 1.644399255026656, 1.660366772655919, 1.676671430718741, 1.685096899472327, 1.688433591127435, 1.687708456444994, 1.685368475440569, 1.683954236709752, 1.676141195382414, 1.671607230319196, 1.669094259348498, 1.668259665855501, 1.672080054237043, 1.686051808929903, 1.69603557244983, 1.690628671719764, 1.671865998881486, 1.655248938188578)
for ($i = 0; $i -lt $colI.Length; $i++) {
    $ws.Cells.Item($i + 2, 9).Value = $colI[$i]
}

$colL = @(0.2236549557337781, 0.2209323942782646, 0.2193779890623802, 0.2187740025749392, 0.2186754884493212, 0.2193697243128128, 0.2226918539530942, 0.2301394042953859, 0.2361842465866033, 0.2390596994956837, 0.2401666930447846, 0.2399274752625331, 0.2391504089941634, 0.238676795219277, 0.2359988617049709, 0.2343882516543943, 0.2334736921001337, 0.2331660664606261, 0.23455848003303, 0.2393781599011788, 0.2426337560237357, 0.2408864959769375, 0.2344814842739424, 0.2280243248771399)
for ($i = 0; $i -lt $colL.Length; $i++) {
    $ws.Cells.Item($i + 2, 12).Value = $colL[$i]
}

$colM = @(0.4085453204198757, 0.3928655311267661, 0.3834853895524191, 0.379724856078326, 0.3791041571795191, 0.3834344231611269, 0.4030874483507532, 0.4436027735483208, 0.4745965229686036, 0.4889680736322788, 0.4944497260098402, 0.4932673943815544, 0.489418259949467, 0.4870657014971442, 0.4736628136079233, 0.4655105365639969, 0.460847178028942, 0.4592726355695405, 0.4663757071364429, 0.4905477707563435, 0.5065757466058045, 0.4980001569240144, 0.4659844903638586, 0.4324284692619855)
for ($i = 0; $i -lt $colM.Length; $i++) {
    $ws.Cells.Item($i + 2, 13).Value = $colM[$i]
}

$colN = @(2.099581378194465, 2.120333747197954, 2.133803042109761, 2.139474297073356, 2.140427008079335, 2.133878788716046, 2.106585572930591, 2.05885443934195, 2.027343705917382, 2.013787046573427, 2.008765873624483, 2.00984226255455, 2.013371694704965, 2.015548234843301, 2.028245366733557, 2.036234267253846, 2.040902438637936, 2.042495549557877, 2.035376258306378, 2.012331958488225, 1.997926848293062, 2.005554930781052, 2.035763929822423, 2.071144266104746)
for ($i = 0; $i -lt $colN.Length; $i++) {
    $ws.Cells.Item($i + 2, 14).Value = $colN[$i]
}
